# FAST_holdings.xlsx update: refresh the "as of" date in the confidential
# disclaimer and refresh the Weight / Percent Change figures for the model
# holdings table (rows 2-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; temporarily unprotect so the cell values can be
# updated, then reapply protection once the edits are done.
$ws.Unprotect()

# --- Weight (column D) and Percent Change (column E) updates ---
$ws.Range("D2").Value2 = 0.08983430837443938
$ws.Range("E2").Value2 = -0.002220480668756686

$ws.Range("D3").Value2 = 0.1056521909210094
$ws.Range("E3").Value2 = -0.006405563689604543

$ws.Range("D4").Value2 = 0.120549379890317
$ws.Range("E4").Value2 = -0.004454168945846604

$ws.Range("D5").Value2 = 0.1417230184516677
$ws.Range("E5").Value2 = -0.001721512964871996

$ws.Range("D6").Value2 = 0.1380561096216016
$ws.Range("E6").Value2 = -0.00399229074889873

$ws.Range("D7").Value2 = 0.1476381978660237
$ws.Range("E7").Value2 = 0.00216267042783258

$ws.Range("D8").Value2 = 0.126424730341808
$ws.Range("E8").Value2 = -0.007516536380036132

$ws.Range("D9").Value2 = 0.1301220645331332
$ws.Range("E9").Value2 = -0.006244067933520392

# Row 10 is the "Total" row; only the Percent Change value changes.
$ws.Range("E10").Value2 = -0.003651796961334441

# --- Update the "as of" date inside the confidential disclaimer text ---
$disclaimer = $ws.Range("A13").Value2
$disclaimer = $disclaimer -replace "2021-05-14", "2021-05-17"
$ws.Range("A13").Value2 = $disclaimer

# Restore the row height (the text length didn't change, but some engines
# auto-adjust row height when cell text is rewritten).
$ws.Rows.Item(13).AutoFit()

# Reinstate sheet protection.
$ws.Protect()
